$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold values in rows 2-4
$ws.Range("B2").Value = 5.3
$ws.Range("B3").Value = 5.4
$ws.Range("B4").Value = 0.8
$ws.Range("C4").Value = 1.4

# Delete row 5 (theta_threshold_range) entirely, shifting row 6 (pie_threshold_range) up to row 5
$ws.Rows(5).Delete()

# Update the (now) row 5 values - pie_threshold_range
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Update selection to reflect the new active cell
$ws.Range("C5").Select() | Out-Null

# Configure page setup for printing (PaperSize 9 = A4, Orientation 1 = portrait)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

